# Updates cryptos list values per latest data refresh (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @{
    "D2" = "37.433.50"
    "E2" = "  -0.87%  "
    "D3" = "2.063.96"
    "E3" = "  -1.06%  "
    "E4" = "  +0.01%  "
    "D5" = "231.73"
    "E5" = "  -0.78%  "
    "D6" = "0.621"
    "E6" = "  -0.77%  "
    "E7" = "  +0.11%  "
    "D8" = "57.66"
    "E8" = "  -2.05%  "
    "D9" = "0.388"
    "E10" = "  -1.16%  "
    "E11" = "  +0.09%  "
    "D12" = "14.79"
    "E12" = "  +0.36%  "
    "D13" = "2.371.22"
    "E13" = "  -0.91%  "
    "D14" = "21.06"
    "E14" = "  -0.96%  "
    "D15" = "0.762"
    "E15" = "  -2.05%  "
    "D16" = "5.31"
    "E16" = "  -0.63%  "
    "D17" = "2.064.76"
    "E17" = "  -0.70%  "
    "D18" = "37.363.17"
    "E18" = "  -1.00%  "
    "D19" = "6.13"
    "E19" = "  -0.89%  "
    "D20" = "70.17"
    "E20" = "  -2.33%  "
    "E21" = "  -2.26%  "
    "D22" = "227.21"
    "E22" = "  -0.48%  "
    "E23" = "  +0.03%  "
    "E24" = "  -0.03%  "
    "D25" = "2.33"
    "E25" = "  -3.29%  "
    "E26" = "  +3.59%  "
    "D27" = "169.31"
    "D28" = "0.130"
    "E28" = "  -5.37%  "
    "D29" = "19.27"
    "E29" = "  -1.55%  "
    "E30" = "  -4.39%  "
    "E31" = "  -0.42%  "
    "E32" = "  -3.49%  "
    "D33" = "0.0629"
    "E33" = "  -1.04%  "
    "D34" = "4.66"
    "E34" = "  -0.57%  "
    "E35" = "  +0.39%  "
    "E36" = "  -0.25%  "
    "D37" = "3.30"
    "E37" = "  -3.95%  "
    "E38" = "  +0.01%  "
    "D39" = "5.30"
    "E39" = "  -2.46%  "
    "D40" = "0.0226"
    "E40" = "  +2.94%  "
    "B41" = "HuobiToken"
    "C41" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D41" = "2.92"
    "E41" = "  +0.01%  "
    "B42" = "Aave"
    "C42" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D42" = "98.05"
    "E42" = "  -1.11%  "
    "D43" = "1.489.53"
    "E43" = "  +2.61%  "
    "D44" = "0.0956"
    "E44" = "  -2.67%  "
    "D45" = "16.96"
    "E45" = "  -1.16%  "
    "D46" = "1.18"
    "E46" = "  +2.49%  "
    "D47" = "4.05"
    "E47" = "  -2.85%  "
    "E48" = "  -2.58%  "
    "D49" = "7.25"
    "E49" = "  -1.62%  "
    "D50" = "2.95"
    "E50" = "  -1.40%  "
    "D51" = "2.256.79"
    "E51" = "  -0.93%  "
}

foreach ($ref in $cellUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$ref]
}
